$d = $word.ActiveDocument

# --- Helper: wrap a document.xml fragment in the pkg:package envelope that
# Range.InsertXML() expects, then apply it to $range (fully replacing the
# range's contents with the supplied paragraph XML). ---
function Set-ParagraphXml($range, $paragraphXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $paragraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# --- Helper: find the (first) paragraph whose text, ignoring the trailing
# paragraph mark, equals $text. ---
function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs($i)
        if ($cand.Range.Text.TrimEnd("`r") -eq $text) {
            return $cand
        }
    }
    return $null
}

# 1) Title paragraph: "Documento de Diseño de " + "Wireframes" (two runs,
#    separated by proofErr spell-check markers) -> single merged run.
$titlePara = Find-ParagraphByText $d "Documento de Diseño de Wireframes"
if ($null -eq $titlePara) { $titlePara = $d.Paragraphs(7) }
$titleXml = '<w:p w14:paraId="00000007" w14:textId="77777777" w:rsidR="003D385A" w:rsidRDefault="00000000"><w:pPr><w:spacing w:before="240" w:after="240"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Documento de Diseño de Wireframes</w:t></w:r></w:p>'
Set-ParagraphXml $titlePara.Range $titleXml

# 2) Author paragraph: "Mendoza " + "Gomez" + ", Carlos Daniel" (three runs,
#    proofErr markers around "Gomez") -> single merged run.
$authorPara = Find-ParagraphByText $d "Mendoza Gomez, Carlos Daniel"
if ($null -eq $authorPara) { $authorPara = $d.Paragraphs(8) }
$authorXml = '<w:p w14:paraId="00000008" w14:textId="77777777" w:rsidR="003D385A" w:rsidRDefault="00000000"><w:pPr><w:spacing w:before="240" w:after="240"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Mendoza Gomez, Carlos Daniel</w:t></w:r></w:p>'
Set-ParagraphXml $authorPara.Range $authorXml

# 3) "CU005: " use-case paragraph: append a new run with "Realizar encuesta".
$cu005Para = Find-ParagraphByText $d "CU005: "
if ($null -ne $cu005Para) {
    $cu005Xml = '<w:p w14:paraId="07A42C60" w14:textId="20C7AA03" w:rsidR="00327A06" w:rsidRDefault="00327A06" w:rsidP="00327A06"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">CU005: </w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Realizar encuesta</w:t></w:r></w:p>'
    Set-ParagraphXml $cu005Para.Range $cu005Xml
} else {
    Write-Host "WARNING: CU005 paragraph not found"
}

Write-Host "Done"
